# Apply the "removed old user stories" edit to the
# StoriesDetailsIteration1 sheet: mark several rows as "Done" in column K,
# widen column J a bit, and update the scroll/selection/zoom state of the
# sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StoriesDetailsIteration1")

# Make sure we're working on the right sheet.
$ws.Activate()

# Mark rows as "Done" in column K.
$doneRows = 5, 8, 13, 19, 26, 29, 30
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 11).Value = "Done"
}

# Widen column J (column 10) - originally unset/default, now custom width.
$ws.Columns.Item(10).ColumnWidth = 20.4

# Update the view: scroll position, zoom, and selection.
$excel.ActiveWindow.Zoom = 70
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L27").Select()
